$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot current (pre-edit) values for the columns that move, rows 2-17
$snapshot = @{}
$snapshot["D2"] = $ws.Range("D2").Value2
$snapshot["M2"] = $ws.Range("M2").Value2
$snapshot["N2"] = $ws.Range("N2").Value2
$snapshot["O2"] = $ws.Range("O2").Value2
$snapshot["P2"] = $ws.Range("P2").Value2
$snapshot["R2"] = $ws.Range("R2").Value2
$snapshot["S2"] = $ws.Range("S2").Value2
$snapshot["D3"] = $ws.Range("D3").Value2
$snapshot["M3"] = $ws.Range("M3").Value2
$snapshot["N3"] = $ws.Range("N3").Value2
$snapshot["O3"] = $ws.Range("O3").Value2
$snapshot["P3"] = $ws.Range("P3").Value2
$snapshot["R3"] = $ws.Range("R3").Value2
$snapshot["S3"] = $ws.Range("S3").Value2
$snapshot["D4"] = $ws.Range("D4").Value2
$snapshot["M4"] = $ws.Range("M4").Value2
$snapshot["N4"] = $ws.Range("N4").Value2
$snapshot["O4"] = $ws.Range("O4").Value2
$snapshot["P4"] = $ws.Range("P4").Value2
$snapshot["R4"] = $ws.Range("R4").Value2
$snapshot["S4"] = $ws.Range("S4").Value2
$snapshot["D5"] = $ws.Range("D5").Value2
$snapshot["M5"] = $ws.Range("M5").Value2
$snapshot["N5"] = $ws.Range("N5").Value2
$snapshot["O5"] = $ws.Range("O5").Value2
$snapshot["P5"] = $ws.Range("P5").Value2
$snapshot["R5"] = $ws.Range("R5").Value2
$snapshot["S5"] = $ws.Range("S5").Value2
$snapshot["D6"] = $ws.Range("D6").Value2
$snapshot["M6"] = $ws.Range("M6").Value2
$snapshot["N6"] = $ws.Range("N6").Value2
$snapshot["O6"] = $ws.Range("O6").Value2
$snapshot["P6"] = $ws.Range("P6").Value2
$snapshot["R6"] = $ws.Range("R6").Value2
$snapshot["S6"] = $ws.Range("S6").Value2
$snapshot["D7"] = $ws.Range("D7").Value2
$snapshot["M7"] = $ws.Range("M7").Value2
$snapshot["N7"] = $ws.Range("N7").Value2
$snapshot["O7"] = $ws.Range("O7").Value2
$snapshot["P7"] = $ws.Range("P7").Value2
$snapshot["R7"] = $ws.Range("R7").Value2
$snapshot["S7"] = $ws.Range("S7").Value2
$snapshot["D8"] = $ws.Range("D8").Value2
$snapshot["M8"] = $ws.Range("M8").Value2
$snapshot["N8"] = $ws.Range("N8").Value2
$snapshot["O8"] = $ws.Range("O8").Value2
$snapshot["P8"] = $ws.Range("P8").Value2
$snapshot["R8"] = $ws.Range("R8").Value2
$snapshot["S8"] = $ws.Range("S8").Value2
$snapshot["D9"] = $ws.Range("D9").Value2
$snapshot["M9"] = $ws.Range("M9").Value2
$snapshot["N9"] = $ws.Range("N9").Value2
$snapshot["O9"] = $ws.Range("O9").Value2
$snapshot["P9"] = $ws.Range("P9").Value2
$snapshot["R9"] = $ws.Range("R9").Value2
$snapshot["S9"] = $ws.Range("S9").Value2
$snapshot["D10"] = $ws.Range("D10").Value2
$snapshot["M10"] = $ws.Range("M10").Value2
$snapshot["N10"] = $ws.Range("N10").Value2
$snapshot["O10"] = $ws.Range("O10").Value2
$snapshot["P10"] = $ws.Range("P10").Value2
$snapshot["R10"] = $ws.Range("R10").Value2
$snapshot["S10"] = $ws.Range("S10").Value2
$snapshot["D11"] = $ws.Range("D11").Value2
$snapshot["M11"] = $ws.Range("M11").Value2
$snapshot["N11"] = $ws.Range("N11").Value2
$snapshot["O11"] = $ws.Range("O11").Value2
$snapshot["P11"] = $ws.Range("P11").Value2
$snapshot["R11"] = $ws.Range("R11").Value2
$snapshot["S11"] = $ws.Range("S11").Value2
$snapshot["D12"] = $ws.Range("D12").Value2
$snapshot["M12"] = $ws.Range("M12").Value2
$snapshot["N12"] = $ws.Range("N12").Value2
$snapshot["O12"] = $ws.Range("O12").Value2
$snapshot["P12"] = $ws.Range("P12").Value2
$snapshot["R12"] = $ws.Range("R12").Value2
$snapshot["S12"] = $ws.Range("S12").Value2
$snapshot["D13"] = $ws.Range("D13").Value2
$snapshot["M13"] = $ws.Range("M13").Value2
$snapshot["N13"] = $ws.Range("N13").Value2
$snapshot["O13"] = $ws.Range("O13").Value2
$snapshot["P13"] = $ws.Range("P13").Value2
$snapshot["R13"] = $ws.Range("R13").Value2
$snapshot["S13"] = $ws.Range("S13").Value2
$snapshot["D14"] = $ws.Range("D14").Value2
$snapshot["M14"] = $ws.Range("M14").Value2
$snapshot["N14"] = $ws.Range("N14").Value2
$snapshot["O14"] = $ws.Range("O14").Value2
$snapshot["P14"] = $ws.Range("P14").Value2
$snapshot["R14"] = $ws.Range("R14").Value2
$snapshot["S14"] = $ws.Range("S14").Value2
$snapshot["D15"] = $ws.Range("D15").Value2
$snapshot["M15"] = $ws.Range("M15").Value2
$snapshot["N15"] = $ws.Range("N15").Value2
$snapshot["O15"] = $ws.Range("O15").Value2
$snapshot["P15"] = $ws.Range("P15").Value2
$snapshot["R15"] = $ws.Range("R15").Value2
$snapshot["S15"] = $ws.Range("S15").Value2
$snapshot["D16"] = $ws.Range("D16").Value2
$snapshot["M16"] = $ws.Range("M16").Value2
$snapshot["N16"] = $ws.Range("N16").Value2
$snapshot["O16"] = $ws.Range("O16").Value2
$snapshot["P16"] = $ws.Range("P16").Value2
$snapshot["R16"] = $ws.Range("R16").Value2
$snapshot["S16"] = $ws.Range("S16").Value2
$snapshot["D17"] = $ws.Range("D17").Value2
$snapshot["M17"] = $ws.Range("M17").Value2
$snapshot["N17"] = $ws.Range("N17").Value2
$snapshot["O17"] = $ws.Range("O17").Value2
$snapshot["P17"] = $ws.Range("P17").Value2
$snapshot["R17"] = $ws.Range("R17").Value2
$snapshot["S17"] = $ws.Range("S17").Value2

# Re-distribute rows according to the target permutation (dest row = row getting the data; src row = row the data came from)
$ws.Range("D2").Value = $snapshot["D6"]
$ws.Range("M2").Value = $snapshot["M6"]
$ws.Range("N2").Value = $snapshot["N6"]
$ws.Range("O2").Value = $snapshot["O6"]
$ws.Range("P2").Value = $snapshot["P6"]
$ws.Range("R2").Value = $snapshot["R6"]
$ws.Range("S2").Value = $snapshot["S6"]
$ws.Range("D3").Value = $snapshot["D2"]
$ws.Range("M3").Value = $snapshot["M2"]
$ws.Range("N3").Value = $snapshot["N2"]
$ws.Range("O3").Value = $snapshot["O2"]
$ws.Range("P3").Value = $snapshot["P2"]
$ws.Range("R3").Value = $snapshot["R2"]
$ws.Range("S3").Value = $snapshot["S2"]
$ws.Range("D4").Value = $snapshot["D5"]
$ws.Range("M4").Value = $snapshot["M5"]
$ws.Range("N4").Value = $snapshot["N5"]
$ws.Range("O4").Value = $snapshot["O5"]
$ws.Range("P4").Value = $snapshot["P5"]
$ws.Range("R4").Value = $snapshot["R5"]
$ws.Range("S4").Value = $snapshot["S5"]
$ws.Range("D5").Value = $snapshot["D17"]
$ws.Range("M5").Value = $snapshot["M17"]
$ws.Range("N5").Value = $snapshot["N17"]
$ws.Range("O5").Value = $snapshot["O17"]
$ws.Range("P5").Value = $snapshot["P17"]
$ws.Range("R5").Value = $snapshot["R17"]
$ws.Range("S5").Value = $snapshot["S17"]
$ws.Range("D6").Value = $snapshot["D13"]
$ws.Range("M6").Value = $snapshot["M13"]
$ws.Range("N6").Value = $snapshot["N13"]
$ws.Range("O6").Value = $snapshot["O13"]
$ws.Range("P6").Value = $snapshot["P13"]
$ws.Range("R6").Value = $snapshot["R13"]
$ws.Range("S6").Value = $snapshot["S13"]
$ws.Range("D8").Value = $snapshot["D14"]
$ws.Range("M8").Value = $snapshot["M14"]
$ws.Range("N8").Value = $snapshot["N14"]
$ws.Range("O8").Value = $snapshot["O14"]
$ws.Range("P8").Value = $snapshot["P14"]
$ws.Range("R8").Value = $snapshot["R14"]
$ws.Range("S8").Value = $snapshot["S14"]
$ws.Range("D9").Value = $snapshot["D11"]
$ws.Range("M9").Value = $snapshot["M11"]
$ws.Range("N9").Value = $snapshot["N11"]
$ws.Range("O9").Value = $snapshot["O11"]
$ws.Range("P9").Value = $snapshot["P11"]
$ws.Range("R9").Value = $snapshot["R11"]
$ws.Range("S9").Value = $snapshot["S11"]
$ws.Range("D10").Value = $snapshot["D4"]
$ws.Range("M10").Value = $snapshot["M4"]
$ws.Range("N10").Value = $snapshot["N4"]
$ws.Range("O10").Value = $snapshot["O4"]
$ws.Range("P10").Value = $snapshot["P4"]
$ws.Range("R10").Value = $snapshot["R4"]
$ws.Range("S10").Value = $snapshot["S4"]
$ws.Range("D11").Value = $snapshot["D9"]
$ws.Range("M11").Value = $snapshot["M9"]
$ws.Range("N11").Value = $snapshot["N9"]
$ws.Range("O11").Value = $snapshot["O9"]
$ws.Range("P11").Value = $snapshot["P9"]
$ws.Range("R11").Value = $snapshot["R9"]
$ws.Range("S11").Value = $snapshot["S9"]
$ws.Range("D12").Value = $snapshot["D15"]
$ws.Range("M12").Value = $snapshot["M15"]
$ws.Range("N12").Value = $snapshot["N15"]
$ws.Range("O12").Value = $snapshot["O15"]
$ws.Range("P12").Value = $snapshot["P15"]
$ws.Range("R12").Value = $snapshot["R15"]
$ws.Range("S12").Value = $snapshot["S15"]
$ws.Range("D13").Value = $snapshot["D16"]
$ws.Range("M13").Value = $snapshot["M16"]
$ws.Range("N13").Value = $snapshot["N16"]
$ws.Range("O13").Value = $snapshot["O16"]
$ws.Range("P13").Value = $snapshot["P16"]
$ws.Range("R13").Value = $snapshot["R16"]
$ws.Range("S13").Value = $snapshot["S16"]
$ws.Range("D14").Value = $snapshot["D12"]
$ws.Range("M14").Value = $snapshot["M12"]
$ws.Range("N14").Value = $snapshot["N12"]
$ws.Range("O14").Value = $snapshot["O12"]
$ws.Range("P14").Value = $snapshot["P12"]
$ws.Range("R14").Value = $snapshot["R12"]
$ws.Range("S14").Value = $snapshot["S12"]
$ws.Range("D15").Value = $snapshot["D10"]
$ws.Range("M15").Value = $snapshot["M10"]
$ws.Range("N15").Value = $snapshot["N10"]
$ws.Range("O15").Value = $snapshot["O10"]
$ws.Range("P15").Value = $snapshot["P10"]
$ws.Range("R15").Value = $snapshot["R10"]
$ws.Range("S15").Value = $snapshot["S10"]
$ws.Range("D16").Value = $snapshot["D3"]
$ws.Range("M16").Value = $snapshot["M3"]
$ws.Range("N16").Value = $snapshot["N3"]
$ws.Range("O16").Value = $snapshot["O3"]
$ws.Range("P16").Value = $snapshot["P3"]
$ws.Range("R16").Value = $snapshot["R3"]
$ws.Range("S16").Value = $snapshot["S3"]
$ws.Range("D17").Value = $snapshot["D8"]
$ws.Range("M17").Value = $snapshot["M8"]
$ws.Range("N17").Value = $snapshot["N8"]
$ws.Range("O17").Value = $snapshot["O8"]
$ws.Range("P17").Value = $snapshot["P8"]
$ws.Range("R17").Value = $snapshot["R8"]
$ws.Range("S17").Value = $snapshot["S8"]
